$p = $ppt.ActivePresentation

# --- 1. Update the cached "today" text on every Date placeholder (master +
#        all slide layouts) from 04/09/2020 to 16/09/2020 --------------------
$newDate = "16/09/2020"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

Update-DateShapes $p.SlideMaster.Shapes
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# --- 2. Remove the old opening "Contexto de Negócio" title slide, the old
#        "Controle da Auto Peças" slide, and the trailing "Resumo" slide -----
$p.Slides.Item(7).Delete()
$p.Slides.Item(2).Delete()
$p.Slides.Item(1).Delete()
